$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# The title placeholder ("Title 1") holds the combined paper title. It
# currently repeats "The Google File System" a second time by mistake,
# right after "... Large-Scale Data Analysis". Fix the typo by removing
# that duplicated trailing text while leaving everything else in the
# run (and the rest of the title) untouched.
$sh = $s.Shapes.Item(1)
$tr = $sh.TextFrame.TextRange

$badTail = "Analysis The Google File System  "
$goodTail = "Analysis "

$fullText = $tr.Text
$startPos = $fullText.IndexOf($badTail)

if ($startPos -ge 0) {
    $target = $tr.Characters($startPos + 1, $badTail.Length)
    $target.Text = $goodTail
}
